$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: B=2 (Tipo Doc), C=3 (N Doc), D=4 (Nombre), E=5 (Periodo Mora), F=6 (Valor Mora), G=7 (Salario Basico)
$data = @(
    @(16, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1903", 26041, 737717),
    @(17, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1902", 31249, 737717),
    @(18, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1901", 31249, 737717),
    @(19, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1812", 31249, 737717),
    @(20, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1811", 31249, 737717),
    @(21, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1810", 31249, 737717),
    @(22, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1809", 31249, 737717),
    @(23, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1808", 29509, 737717),
    @(24, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1803", 29509, 737717),
    @(25, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1802", 29509, 737717),
    @(26, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1801", 29509, 737717),
    @(27, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1712", 29509, 737717),
    @(28, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1711", 29509, 737717),
    @(29, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1710", 29509, 737717),
    @(30, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1709", 29509, 737717),
    @(31, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1708", 29509, 737717),
    @(32, "CC", "7918105", "PEDRO EDWIN BELLO BERRIO", "1707", 29509, 737717),
    @(33, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1805", 31249, 781242),
    @(34, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1803", 29509, 781242),
    @(35, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1802", 29509, 781242),
    @(36, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1711", 29509, 781242),
    @(37, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1710", 29509, 781242),
    @(38, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1709", 29509, 781242),
    @(39, "CC", "9098844", "VICTOR ISAAC MURILLO JULIO", "1708", 29509, 781242),
    @(40, "CC", "8851958", "RONALD DAVID ARANDA GOMEZ", "1805", 31249, 781242),
    @(41, "CC", "8851958", "RONALD DAVID ARANDA GOMEZ", "1803", 29509, 781242),
    @(42, "CC", "8851958", "RONALD DAVID ARANDA GOMEZ", "1802", 29509, 781242),
    @(43, "CC", "8851958", "RONALD DAVID ARANDA GOMEZ", "1801", 29509, 781242)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

Write-Output "Done updating rows 16-43"